$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing covid_deaths (column C) values per diff
$updates = @(
  @(774, 13),
  @(1147, 3),
  @(1162, 46),
  @(1166, 36),
  @(1170, 8),
  @(1213, 52),
  @(1217, 15),
  @(1218, 28),
  @(1223, 30),
  @(1232, 6),
  @(1240, 35),
  @(1248, 5),
  @(1254, 31),
  @(1260, 35),
  @(1269, 26),
  @(1277, 18),
  @(1282, 36),
  @(1291, 41),
  @(1301, 26),
  @(1306, 17),
  @(1311, 28),
  @(1315, 10),
  @(1316, 22),
  @(1317, 3),
  @(1318, 6),
  @(1319, 12),
  @(1320, 24),
  @(1324, 10),
  @(1325, 16),
  @(1333, 7),
  @(1334, 16),
  @(1335, 20),
  @(1337, 14),
  @(1338, 19),
  @(1341, 8),
  @(1342, 14),
  @(1344, 4),
  @(1346, 15),
  @(1347, 12),
  @(1352, 7),
  @(1357, 12),
  @(1358, 12),
  @(1360, 3),
  @(1363, 12),
  @(1364, 13),
  @(1365, 3),
  @(1366, 9),
  @(1367, 7),
  @(1368, 17)
)
foreach ($u in $updates) {
  $row = $u[0]
  $val = $u[1]
  $ws.Cells.Item($row, 3).Value = $val
}

# Append new rows 1369-1375
$newRows = @(
  @(1369, 44231, "20-29", 1),
  @(1370, 44231, "50-59", 1),
  @(1371, 44231, "60-69", 5),
  @(1372, 44231, "70-79", 4),
  @(1373, 44231, "80+", 16),
  @(1374, 44232, "70-79", 2),
  @(1375, 44232, "80+", 5)
)
foreach ($r in $newRows) {
  $row = $r[0]
  $dateVal = $r[1]
  $agegrp = $r[2]
  $deaths = $r[3]
  $ws.Cells.Item($row, 1).Value = $dateVal
  $ws.Cells.Item($row, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
  $ws.Cells.Item($row, 2).Value = $agegrp
  $ws.Cells.Item($row, 3).Value = $deaths
}
